# Fix several bugs for DcvAdmission Permit
#
# 1) Insert a new "Departure Date" paragraph (mirroring the existing
#    "Arrival Date" paragraph) right after the "Arrival Date" paragraph.
# 2) Collapse the "Water based:" run-soup for arrival.adult.water_based
#    down to a single run.
# 3) Collapse the "Water based:" run-soup for arrival.child.water_based
#    down to two runs (matching the already-split "{" / "{ ... }}" style
#    used by the "Extended stay:" paragraph above it).

$d = $word.ActiveDocument

function Get-ParagraphByText($needle) {
    $rng = $d.Content
    [void]$rng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $rng.Paragraphs(1)
}

# ---------------------------------------------------------------------
# 1) Insert the new "Departure Date" paragraph after "Arrival Date".
# ---------------------------------------------------------------------
$arrivalPara = Get-ParagraphByText("Arrival Date")
$insertPos = $arrivalPara.Range.End - 1   # just before the pilcrow, so we add a *new* paragraph
$insertRange = $d.Range($insertPos, $insertPos)

$departureXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Normal"/>
    <w:tabs>
      <w:tab w:val="clear" w:pos="720"/>
      <w:tab w:val="left" w:pos="3261" w:leader="none"/>
    </w:tabs>
    <w:spacing w:before="0" w:after="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial Nova" w:hAnsi="Arial Nova" w:cs="Arial"/>
      <w:color w:val="464646"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Arial" w:ascii="Arial" w:hAnsi="Arial"/>
      <w:color w:val="464646"/>
    </w:rPr>
    <w:t>Departure</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Arial" w:ascii="Arial" w:hAnsi="Arial"/>
      <w:color w:val="464646"/>
    </w:rPr>
    <w:t xml:space="preserve"> Date</w:t>
    <w:tab/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Arial" w:ascii="Arial Nova" w:hAnsi="Arial Nova"/>
      <w:color w:val="464646"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
    </w:rPr>
    <w:t>{{ arrival.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Arial" w:ascii="Arial Nova" w:hAnsi="Arial Nova"/>
      <w:color w:val="464646"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
    </w:rPr>
    <w:t>departure</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Arial" w:ascii="Arial Nova" w:hAnsi="Arial Nova"/>
      <w:color w:val="464646"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
    </w:rPr>
    <w:t>_date }}</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
[void]$insertRange.InsertXML($departureXml)

# The OM round-trip drops an explicit w:before="0" (it reads as the
# paragraph-format default) - force it back so the new paragraph's
# <w:spacing> matches "Arrival Date"'s sibling exactly.
$departurePara = Get-ParagraphByText("Departure")
$departurePara.Format.SpaceBefore = 0
$departurePara.Format.SpaceAfter = 0

# ---------------------------------------------------------------------
# 2) Simplify the adult "Water based:" paragraph to a single run.
# ---------------------------------------------------------------------
$adultPara = Get-ParagraphByText("{{ arrival.adult.water_based }}")
$adultRange = $adultPara.Range
$adultRange.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
$adultPos = $adultRange.Start
[void]$adultRange.Delete()
# Re-seat a fresh Range at the same spot - reusing the post-Delete Range
# object directly for InsertXML loses the paragraph's own <w:pPr>.
$adultRange = $d.Range($adultPos, $adultPos)

$adultXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Arial" w:ascii="Arial Nova Light" w:hAnsi="Arial Nova Light"/>
      <w:color w:val="464646"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
    </w:rPr>
    <w:t>Water based:</w:t>
    <w:tab/>
    <w:t>{{ arrival.adult.water_based }}</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
[void]$adultRange.InsertXML($adultXml)

# ---------------------------------------------------------------------
# 3) Simplify the child "Water based:" paragraph to two runs.
# ---------------------------------------------------------------------
$childPara = Get-ParagraphByText("{ arrival.child.water_based }}")
$childRange = $childPara.Range
$childRange.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
$childPos = $childRange.Start
[void]$childRange.Delete()
# Re-seat a fresh Range at the same spot - reusing the post-Delete Range
# object directly for InsertXML loses the paragraph's own <w:pPr>.
$childRange = $d.Range($childPos, $childPos)

$childXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Arial" w:ascii="Arial Nova Light" w:hAnsi="Arial Nova Light"/>
      <w:color w:val="464646"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
    </w:rPr>
    <w:t>Water based:</w:t>
    <w:tab/>
    <w:t>{</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Arial" w:ascii="Arial Nova" w:hAnsi="Arial Nova"/>
      <w:color w:val="464646"/>
      <w:sz w:val="18"/>
      <w:szCs w:val="18"/>
    </w:rPr>
    <w:t>{ arrival.child.water_based }}</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
[void]$childRange.InsertXML($childXml)

Write-Output "ok"
